$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fill in the new modification row (row 12) that was previously blank
$ws.Range("A12").Value = "GA03"
$ws.Range("B12").Value = "Gabriel Simard"
$ws.Range("D12").Value = "NomVersion de varchar(35) à varchar(50)"
$ws.Range("E12").Value = "Terminée"

# Match the date formatting already used by the rows above (e.g. C10/C11)
$ws.Range("C10").Copy()
$ws.Range("C12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C12").Value = 42262

# Update the active selection to reflect where the user last worked
$ws.Range("D12").Select()
